# Add the new "Sprachangaben" worksheet (refs #120, refs #134) and make it
# the active/selected sheet, mirroring the layout of the existing
# "Fingerprint" sheet.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("Fingerprint")

# Insert the new sheet right after the last existing sheet ("Fingerprint").
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Sprachangaben"

# Copy the cell formatting (styles, column-esque look & feel) from the
# "Fingerprint" sheet's A1:G12 block so the new sheet matches the
# established template for these description sheets.
$src.Range("A1:G12").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Column widths / row heights matching the template sheet.
$ws.Columns.Item(1).ColumnWidth = 25.7209302325581
$ws.Columns.Item(2).ColumnWidth = 36.7953488372093
$ws.Columns.Item(3).ColumnWidth = 46.6418604651163
$ws.Columns.Item(4).ColumnWidth = 25.106976744186
$ws.Columns.Item(5).ColumnWidth = 26.9488372093023
$ws.Columns.Item(6).ColumnWidth = 22.7674418604651
$ws.Columns.Item(7).ColumnWidth = 45.7813953488372

$ws.Rows.Item(1).RowHeight = 14.15
$ws.Rows.Item(2).RowHeight = 39.55
$ws.Rows.Item(3).RowHeight = 26.85
$ws.Rows.Item(5).RowHeight = 14.2
$ws.Rows.Item(7).RowHeight = 14.2
$ws.Rows.Item(8).RowHeight = 14.15
$ws.Rows.Item(10).RowHeight = 26.95
$ws.Rows.Item(11).RowHeight = 24.7
$ws.Rows.Item(12).RowHeight = 14.15

# Header block.
$ws.Range("A1").Value = "Titel"
$ws.Range("B1").Value = "Sprachangaben"

$ws.Range("A2").Value = "Sprachliche Beschreibung"
$ws.Range("B2").Value = "Die einzelnen Sprachen (Marc 041a) werden nacheinander angezeigt und durch ein Komma getrennt."

$ws.Range("A3").Value = "Ungefähre Entsprechung  Marc"
$ws.Range("B3").Value = "Language Code"

$ws.Range("A5").Value = "Anforderung in Redmine"
$ws.Range("B5").Value = "Ticket #134"
$ws.Hyperlinks.Add($ws.Range("B5"), "http://redmine.thulb.uni-jena.de/issues/134", "", "", "Ticket #134")

# Datenfelder block.
$ws.Range("A7").Value = "Datenfelder"
$ws.Range("B7").Value = "Bemerkung"
$ws.Range("C7").Value = "Anmerkung (Fachabteilung)"

$ws.Range("A8").Value = "546 `$a-`$3"

# Genutzte Felder block.
$ws.Range("A10").Value = "Test Kommentar"
$ws.Range("B10").Value = "Test Titel (PPN)"
$ws.Range("C10").Value = "Anzuzeigende Information in Vollanzeige"
$ws.Range("D10").Value = "Anzuzeigende Information in der Kurzanzeige"
$ws.Range("E10").Value = "enthaltener Link"

$ws.Range("A11").Value = "Genutzte Felder"

$ws.Range("A12").Value = "041 `$a"
$ws.Range("B12").Value = 786233990
$ws.Range("C12").Value = "Deutsch, Französisch"

# Re-apply the formats for the row/cells touched after the copy/paste so
# the value-setting above doesn't disturb the pasted look (Value alone
# does not change formatting, this is just a safety re-assert).
$src.Range("A1:G12").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# The new sheet becomes the active tab (and the only "selected"/tabSelected
# sheet); previously it was "Fingerprint" (activeTab index 1 -> 4).
$ws.Activate()
